# Commit: "Fruta / hortaliza, semanal"
#
# Three new weekly price records for Choclo (Vega Central Mapocho de
# Santiago) are inserted right before the existing row 715, pushing all
# subsequent rows (old 715-768) down by three (new 718-771).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three blank rows at 715-717; everything below shifts down.
$ws.Rows("715:717").Insert()

# ---- New row 715 ----
$ws.Range("A715").Value = 9
$ws.Range("B715").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C715").Value = "Metropolitana"
$ws.Range("D715").Value = 45021
$ws.Range("E715").Value = 13
$ws.Range("F715").Value = 100112024
$ws.Range("G715").Value = "Choclo"
$ws.Range("H715").Value = "Choclero"
$ws.Range("I715").Value = "Primera"
$ws.Range("J715").Value = 9800
$ws.Range("K715").Value = 300
$ws.Range("L715").Value = 350
$ws.Range("M715").Value = 324
$ws.Range("N715").Value = "`$/unidad"
$ws.Range("O715").Value = "Región Metropolitana"
$ws.Range("P715").Value = 324
$ws.Range("Q715").Value = 1
$ws.Range("R715").Value = "Hortaliza"

# ---- New row 716 ----
$ws.Range("A716").Value = 9
$ws.Range("B716").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C716").Value = "Metropolitana"
$ws.Range("D716").Value = 45021
$ws.Range("E716").Value = 13
$ws.Range("F716").Value = 100112024
$ws.Range("G716").Value = "Choclo"
$ws.Range("H716").Value = "Choclero"
$ws.Range("I716").Value = "Segunda"
$ws.Range("J716").Value = 5500
$ws.Range("K716").Value = 250
$ws.Range("L716").Value = 250
$ws.Range("M716").Value = 250
$ws.Range("N716").Value = "`$/unidad"
$ws.Range("O716").Value = "Región Metropolitana"
$ws.Range("P716").Value = 250
$ws.Range("Q716").Value = 1
$ws.Range("R716").Value = "Hortaliza"

# ---- New row 717 ----
$ws.Range("A717").Value = 9
$ws.Range("B717").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C717").Value = "Metropolitana"
$ws.Range("D717").Value = 45021
$ws.Range("E717").Value = 13
$ws.Range("F717").Value = 100112024
$ws.Range("G717").Value = "Choclo"
$ws.Range("H717").Value = "Dulce o Americano"
$ws.Range("I717").Value = "Primera"
$ws.Range("J717").Value = 7500
$ws.Range("K717").Value = 250
$ws.Range("L717").Value = 280
$ws.Range("M717").Value = 264
$ws.Range("N717").Value = "`$/unidad"
$ws.Range("O717").Value = "Región Metropolitana"
$ws.Range("P717").Value = 264
$ws.Range("Q717").Value = 1
$ws.Range("R717").Value = "Hortaliza"
